$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = 6
$ws.Range("I10").Value = 1.53
$ws.Range("L10").Value = 1.29
$ws.Range("M10").Value = 3.5
$ws.Range("N10").Value = 1.95
$ws.Range("O10").Value = 1.85
$ws.Range("X10").Value = 51
$ws.Range("AB10").Value = 21
$ws.Range("AE10").Value = 6
$ws.Range("AF10").Value = 6.5
$ws.Range("G14").Value = 2.35
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 9.5
$ws.Range("N14").Value = 2.03
$ws.Range("O14").Value = 1.78
$ws.Range("U14").Value = 11
$ws.Range("X14").Value = 19
$ws.Range("Z14").Value = 9.5
$ws.Range("AB14").Value = 13
$ws.Range("AE14").Value = 9.5
$ws.Range("J16").Value = 1.04
$ws.Range("L16").Value = 1.3
$ws.Range("P16").Value = 1.37
$ws.Range("J17").Value = 1.02
$ws.Range("L17").Value = 1.11
$ws.Range("P17").Value = 1.22
$ws.Range("J18").Value = 1.03
$ws.Range("L18").Value = 1.17
$ws.Range("P18").Value = 1.27
$ws.Range("G22").Value = 3.25
$ws.Range("K22").Value = 9.25
$ws.Range("L22").Value = 1.16
$ws.Range("M22").Value = 4.55
$ws.Range("N22").Value = 1.5
$ws.Range("O22").Value = 2.4
$ws.Range("P22").Value = 1.28
$ws.Range("Q22").Value = 3.35
$ws.Range("R22").Value = 1.47
$ws.Range("S22").Value = 2.52
$ws.Range("T22").Value = 15.5
$ws.Range("U22").Value = 22
$ws.Range("V22").Value = 11.5
$ws.Range("W22").Value = 45
$ws.Range("Y22").Value = 23
$ws.Range("Z22").Value = 9.25
$ws.Range("AA22").Value = 7.8
$ws.Range("AB22").Value = 11.25
$ws.Range("AC22").Value = 35
$ws.Range("AD22").Value = 175
$ws.Range("AE22").Value = 11
$ws.Range("AF22").Value = 12
$ws.Range("AG22").Value = 8.5
$ws.Range("AH22").Value = 19.5
$ws.Range("AI22").Value = 14
$ws.Range("AJ22").Value = 18.5
$ws.Range("O23").Value = 1.57
$ws.Range("N24").Value = 2.2
$ws.Range("O24").Value = 1.65
$ws.Range("G25").Value = 4.5
$ws.Range("K25").Value = 21
$ws.Range("L25").Value = 1.14
$ws.Range("M25").Value = 5.5
$ws.Range("P25").Value = 1.25
$ws.Range("Q25").Value = 3.75
$ws.Range("T25").Value = 19
$ws.Range("V25").Value = 15
$ws.Range("Y25").Value = 29
$ws.Range("AD25").Value = 126
$ws.Range("AI25").Value = 12
$ws.Range("G26").Value = 4.75
$ws.Range("H26").Value = 4.75
$ws.Range("I26").Value = 1.57
$ws.Range("J26").Value = 1.01
$ws.Range("K26").Value = 23
$ws.Range("L26").Value = 1.1
$ws.Range("M26").Value = 7
$ws.Range("X26").Value = 34
$ws.Range("AA26").Value = 10
$ws.Range("AE26").Value = 13
$ws.Range("AH26").Value = 13
$ws.Range("G30").Value = 2.45
$ws.Range("H30").Value = 3.25
$ws.Range("I30").Value = 2.65
$ws.Range("L30").Value = 1.32
$ws.Range("M30").Value = 2.85
$ws.Range("N30").Value = 1.93
$ws.Range("O30").Value = 1.7
$ws.Range("P30").Value = 1.42
$ws.Range("Q30").Value = 2.47
$ws.Range("R30").Value = 1.75
$ws.Range("S30").Value = 1.87
$ws.Range("T30").Value = 8
$ws.Range("Y30").Value = 32
$ws.Range("Z30").Value = 9.25
$ws.Range("AA30").Value = 6.3
$ws.Range("AB30").Value = 14.5
$ws.Range("AC30").Value = 70
$ws.Range("AD30").Value = 600
$ws.Range("AE30").Value = 8.25
$ws.Range("AH30").Value = 30
$ws.Range("AI30").Value = 23
$ws.Range("H31").Value = 3.3
$ws.Range("I31").Value = 3.25
$ws.Range("T31").Value = 8.5
$ws.Range("V31").Value = 9
$ws.Range("AA31").Value = 6.5
$ws.Range("AB31").Value = 13
$ws.Range("AI31").Value = 23
$ws.Range("N32").Value = 1.57
$ws.Range("O32").Value = 2.35
$ws.Range("G34").Value = 3.1
$ws.Range("I34").Value = 2.2
$ws.Range("R34").Value = 1.73
$ws.Range("S34").Value = 2
$ws.Range("U34").Value = 17
$ws.Range("W34").Value = 34
$ws.Range("AE34").Value = 8.5
$ws.Range("AF34").Value = 11
$ws.Range("AI34").Value = 17
$ws.Range("N42").Value = 1.9
$ws.Range("O42").Value = 1.9
$ws.Range("G43").Value = 1.7
$ws.Range("H43").Value = 3.9
$ws.Range("I43").Value = 4.5
$ws.Range("N43").Value = 1.65
$ws.Range("O43").Value = 2.2
$ws.Range("U43").Value = 9
$ws.Range("X43").Value = 13
$ws.Range("AA43").Value = 7.5
$ws.Range("AB43").Value = 13
$ws.Range("AE43").Value = 15
$ws.Range("AF43").Value = 26
$ws.Range("AI43").Value = 34
$ws.Range("AJ43").Value = 34
$ws.Range("N45").Value = 1.88
$ws.Range("O45").Value = 1.93
$ws.Range("J47").Value = 1.05
$ws.Range("K47").Value = 11
$ws.Range("G48").Value = 1.33
$ws.Range("H48").Value = 5.3
$ws.Range("M48").Value = 5.4
$ws.Range("N48").Value = 1.38
$ws.Range("O48").Value = 2.82
$ws.Range("P48").Value = 1.2
$ws.Range("Q48").Value = 4.05
$ws.Range("T48").Value = 11
$ws.Range("Z48").Value = 10.25
$ws.Range("AA48").Value = 11.25
$ws.Range("AB48").Value = 17.5
$ws.Range("AG48").Value = 23
$ws.Range("AH48").Value = 150
$ws.Range("L52").Value = 1.5
$ws.Range("M52").Value = 2.5
